$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "https://aaai-kdf2020.github.io/"
$ws.Range("C2").Value = "AAAI Tutorial on Feb 8 2-6pm (Gibson, 2nd floor)"
$ws.Range("C15").Select()
